$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition) - column F holds "想去人数" (want-to-go count)
$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2 = 0
    3 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 2220
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    12 = 0
    13 = 0
    15 = 0
    16 = 33
    17 = 0
    18 = 0
    19 = 0
    20 = 0
    21 = 3986
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    26 = 0
    27 = 0
    29 = 0
    30 = 95
    31 = 0
    32 = 0
    33 = 0
    34 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Cells.Item($row, 6).Value = $sheet1Updates[$row]
}

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Cells.Item(2, 6).Value = 0

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2 = 0
    4 = 0
    5 = 0
    6 = 0
    7 = 0
    8 = 0
    9 = 0
    10 = 0
    11 = 0
    13 = 0
    14 = 0
    15 = 0
    16 = 0
    19 = 0
    20 = 0
    21 = 0
    22 = 0
    23 = 0
    24 = 0
    25 = 0
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = 0
    33 = 0
    35 = 0
    36 = 0
    37 = 0
    38 = 0
    39 = 0
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Cells.Item($row, 6).Value = $sheet4Updates[$row]
}
